# Update the "取得日時" (retrieved datetime) timestamp for each existing
# data row on the "ランサーズ" sheet to reflect the latest scrape run.
# This corresponds to commit: "Append: 2026-01-27 18:44 JST"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-27 18:44:52"

# Data rows are 2 through 14 (row 1 is the header row).
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
